$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Helper: locate a unique piece of text in the document body and replace the
# whole paragraph that contains it with a freshly authored paragraph (used
# whenever a single run needs to be split into several runs with identical
# rendered text/formatting but separate <w:r> boundaries).
# ---------------------------------------------------------------------------
function Set-ParagraphXml($findText, $paragraphInnerXml) {
    $rng = $d.Content
    $ok = $rng.Find.Execute($findText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $ok) {
        throw "Find failed for: $findText"
    }
    $pkg = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
           '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
           '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
           '<pkg:xmlData>' +
           '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
           '<w:body>' + $paragraphInnerXml + '</w:body>' +
           '</w:document>' +
           '</pkg:xmlData></pkg:part></pkg:package>'
    $rng.InsertXML($pkg)
}

# 1) "1/giorno" -> the run "/giorno" is split into "/" and "giorno"
Set-ParagraphXml "/giorno" '<w:p><w:r><w:t>1</w:t></w:r><w:r><w:t>/</w:t></w:r><w:r><w:t>giorno</w:t></w:r></w:p>'

# 2) "Il <<tempo>> provvede ad avviare l'operazione di calcolo e pubblicazione
#    delle statistiche periodiche" -> split into three runs, changing the verb
#    "provvede ad avviare" to "avvia"
Set-ParagraphXml "Il <<tempo>> provvede ad avviare" ('<w:p><w:pPr><w:ind w:left="34"/></w:pPr>' +
    '<w:r><w:t xml:space="preserve">Il &lt;&lt;tempo&gt;&gt; </w:t></w:r>' +
    '<w:r><w:t>avvia</w:t></w:r>' +
    '<w:r><w:t xml:space="preserve"> l’operazione di calcolo e pubblicazione delle statistiche periodiche</w:t></w:r>' +
    '</w:p>')

# 3) "Il sistema pubblica i dati sovrascrivendo i dati riguardanti il periodo
#    di tempo precedente" -> "Il sistema pubblica i dati elaborati" split in
#    two runs
Set-ParagraphXml "Il sistema pubblica i dati sovrascrivendo" ('<w:p>' +
    '<w:r><w:t xml:space="preserve">Il sistema pubblica i </w:t></w:r>' +
    '<w:r><w:t>dati elaborati</w:t></w:r>' +
    '</w:p>')

# 4) Collapse the double space after "ERRORE:" to a single space (same run).
$rng4 = $d.Content
$null = $rng4.Find.Execute("I Scenario/Flusso di eventi di ERRORE:  ", $true, $false, $false, $false, $false, $true, 1, $false, "I Scenario/Flusso di eventi di ERRORE: ", 2)

# 5) Collapse the double space after "Il" to a single space (same run).
$rng5 = $d.Content
$null = $rng5.Find.Execute("Il  sistema invia una notifica di errore riguardante il reperimento dei dati", $true, $false, $false, $false, $false, $true, 1, $false, "Il sistema invia una notifica di errore riguardante il reperimento dei dati", 2)
